# Update template for mass book upload: reorder/insert header columns on
# Φύλλο1 and tweak the window/view state to match the refreshed layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header order (row 1). Column C now holds the "Θεματική" (subject)
# header, shifting the remaining headers one column to the right, with a
# couple of later headers ("Τόπος Έκδοσης" and "Σχόλια") re-ordered too.
$ws.Range("A1").Value = "Κωδικός"
$ws.Range("B1").Value = "Τίτλος*"
$ws.Range("C1").Value = "Θεματική"
$ws.Range("D1").Value = "Συγγραφέας*"
$ws.Range("E1").Value = "Εκδότης"
$ws.Range("F1").Value = "Χρονολογία Έκδοσης"
$ws.Range("G1").Value = "Σχόλια"
$ws.Range("H1").Value = "Τόπος Έκδοσης"
$ws.Range("I1").Value = "Αριθμός Σελίδων"
$ws.Range("J1").Value = "Τρόπος απόκτησης"
$ws.Range("K1").Value = "Χρονολογία απόκτησης"

# Apply the header style (bold/fill/border/center) used for the rest of the
# row to the newly-populated C1 cell as well.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Column widths shift along with the header: everything from E onward moves
# one column right, and the trailing column (L) picks up the width that
# used to belong to K.
$ws.Columns.Item(5).ColumnWidth = 18.85546875
$ws.Columns.Item(6).ColumnWidth = 18.85546875
$ws.Columns.Item(7).ColumnWidth = 15.7109375
$ws.Columns.Item(8).ColumnWidth = 19.28515625
$ws.Columns.Item(9).ColumnWidth = 19.85546875
$ws.Columns.Item(10).ColumnWidth = 24.7109375
$ws.Columns.Item(11).ColumnWidth = 22.28515625
$ws.Columns.Item(12).ColumnWidth = 16.28515625

# Scroll/selection state: the sheet view now starts at column H with H8
# selected as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("H8").Select()

# Workbook window size/position tweak recorded alongside the sheet change.
$excel.ActiveWindow.Top = 3885
$excel.ActiveWindow.Height = 9000
